{"js": "// The edit replaces every occurrence of \"\u0411\u043e\u043e\u0442\u0435\u0441 \u0446\u043e\u043d\u0441\u0442\u0435\u043b\u043b\u0430\u0442\u0438\u043e\u043d\" with\n// \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441\" throughout the document body. This single\n// substitution rule reproduces all five changed spots in the diff:\n//   - \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441 \u0446\u043e\u043d\u0441\u0442\u0435\u043b\u043b\u0430\u0442\u0438\u043e\u043d\" -> \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441\" (x4)\n//   - \"...\u0443\u043d\u0443\u0442\u0430\u0440  \u0411\u043e\u043e\u0442\u0435\u0441 \u0446\u043e\u043d\u0441\u0442\u0435\u043b\u043b\u0430\u0442\u0438\u043e\u043d \u0438...\" -> \"...\u0443\u043d\u0443\u0442\u0430\u0440  \u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441 \u0438...\" (x1)\nconst body = context.document.body;\nconst results = body.search(\"\u0411\u043e\u043e\u0442\u0435\u0441 \u0446\u043e\u043d\u0441\u0442\u0435\u043b\u043b\u0430\u0442\u0438\u043e\u043d\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The edit replaces every occurrence of \"\u0411\u043e\u043e\u0442\u0435\u0441 \u0446\u043e\u043d\u0441\u0442\u0435\u043b\u043b\u0430\u0442\u0438\u043e\u043d\" with\n# \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441\" throughout the document. This single substitution\n# rule reproduces all five changed spots in the diff:\n#   - \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441 \u0446\u043e\u043d\u0441\u0442\u0435\u043b\u043b\u0430\u0442\u0438\u043e\u043d\" -> \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441\" (x4)\n#   - \"...\u0443\u043d\u0443\u0442\u0430\u0440  \u0411\u043e\u043e\u0442\u0435\u0441 \u0446\u043e\u043d\u0441\u0442\u0435\u043b\u043b\u0430\u0442\u0438\u043e\u043d \u0438...\" -> \"...\u0443\u043d\u0443\u0442\u0430\u0440  \u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441 \u0438...\" (x1)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"\u0411\u043e\u043e\u0442\u0435\u0441 \u0446\u043e\u043d\u0441\u0442\u0435\u043b\u043b\u0430\u0442\u0438\u043e\u043d\"\n$find.Replacement.Text = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u0411\u043e\u043e\u0442\u0435\u0441\"\n\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll) -> replace every match in the document.\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
